# Regenerate s_val data to filter save games
# Update cached B/C/D/E/G values for rows 2-7 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 0.6753301551942219; C = 0.3127903958511391;  D = 3.900430680208489;  E = 0.496779210170732; G = 5.385330441424582 }
    3 = @{ B = 1.459612070389937;  C = 1.667794583268128;   D = 0.8054896365839992; E = 0.496779210170732; G = 4.429675500412797 }
    4 = @{ B = 3.230985683306322;  C = 1.667794583268128;   D = 0.8054896365839992; E = 0.496779210170732; G = 6.201049113329182 }
    5 = @{ B = 3.230985683306322;  C = 1.667794583268128;   D = 3.900430680208489;  E = 0.496779210170732; G = 9.295990156953671 }
    6 = @{ B = 3.230985683306322;  C = 1.667794583268128;   D = 26.21740644021617;  E = 0.496779210170732; G = 31.61296591696135 }
    7 = @{ B = 0.127881588408715;  C = 0.04240448674262143; D = 0.8054896365839992; E = 8.660232485948974; G = 9.63600819768431 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
